$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns D:I should pick up the same formatting
# already applied to A1:C1 (bold font, border, center/top alignment).
$ws.Range("A1").Value = "Job_Id"
$ws.Range("B1").Value = "Job_Title"
$ws.Range("C1").Value = "Job_Description"
$ws.Range("D1").Value = "Total_Years_Min_Exp"
$ws.Range("E1").Value = "Total_Years_Max_Exp"
$ws.Range("F1").Value = "Linked_Poster"
$ws.Range("G1").Value = "Linked_Posted"
$ws.Range("H1").Value = "Resume_received"
$ws.Range("I1").Value = "Resume_downloaded"

$ws.Range("A1").Copy()
$ws.Range("D1:I1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data row (row 2)
$ws.Range("A2").Value = "JD_001"
$ws.Range("B2").Value = "Senior Engineer"
$ws.Range("C2").Value = "We are seeking a Software Engineer to build and maintain high-quality software solutions.`nWork with global teams to drive innovation and deliver scalable applications.`nJoin Akkodis and be part of a tech-driven, collaborative environment."
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0

# A2 previously carried a "wrap text" style; the refreshed data row uses
# plain default formatting, so clear it back to Normal.
$ws.Range("A2").Style = "Normal"

# Writing the multi-line description auto-expanded row 2's height; restore
# the sheet's default (non-custom) row height.
$ws.Rows(2).AutoFit()
